$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Reset to factory defaults" -> "Reset to defaults"   (row 23, col B)
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "Reset to defaults"

# ---------------------------------------------------------------------------
# 2) "Reset options and control bindings?" -> "Reset to defaults?" (row 24, col B)
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = "Reset to defaults?"

# ---------------------------------------------------------------------------
# 3) Add "fps limit" / "FPS limit" row right after "vsync" / "V-sync" (row 32)
#    and before "gamma" / "Gamma" (row 33)
# ---------------------------------------------------------------------------
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = "fps limit"
$ws.Range("B33").Value = "FPS limit"

# ---------------------------------------------------------------------------
# 4) Remove "special ability" / "Special ability" row (between "Roll / Kick"
#    and "toggle phone"). It was at row 62 before the insert above shifted
#    everything below row 33 down by one, so it is now row 63.
# ---------------------------------------------------------------------------
$ws.Rows.Item(63).Delete()

# ---------------------------------------------------------------------------
# 5) Remove the conditional formatting rule on column C (and its backing dxf)
# ---------------------------------------------------------------------------
$ws.Cells.FormatConditions.Delete()

# ---------------------------------------------------------------------------
# 6) Update the saved view state: scroll position + selection
#    (scrolls so row 41 is visible, then selects the whole column C - this
#    mirrors the user's click on the column C header before removing the
#    conditional formatting that used to live on that column)
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A41"), $true)
$ws.Range("C:C").Select()
